$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (columns B..Q) - identical across all data rows (2..26)
# Scientific notation avoided (unsupported by parser); plain decimals used
# with enough precision for exact double round-trip.
$newValues = @(
    0.99999691410162661,
    0.99901695851944816,
    0.99999743747974434,
    0.99998391672400289,
    0.99999500368331118,
    0.00000288055067624043,
    0.0009176260715613002,
    0.000001345678104073282,
    0.000003118444667527866,
    0.000002232061385800574,
    0.00009999963621106398,
    0.00169721851163615,
    0.99997531281301288,
    0.001769472617403162,
    67.51505815013246092,
    93.11145047236465189
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
